$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets item 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3302
$ws.Range("F6").Value = 7717
$ws.Range("F8").Value = 709
$ws.Range("F9").Value = 1128
$ws.Range("F10").Value = 1058
$ws.Range("F14").Value = 362
$ws.Range("F15").Value = 6179
$ws.Range("F16").Value = 76
$ws.Range("F21").Value = 5102
$ws.Range("F22").Value = 5925
$ws.Range("F23").Value = 360
$ws.Range("F24").Value = 160
$ws.Range("F30").Value = 8
$ws.Range("F31").Value = 1052
$ws.Range("F32").Value = 1028
$ws.Range("F41").Value = 411
$ws.Range("F43").Value = 1196
$ws.Range("F47").Value = 3202
$ws.Range("F48").Value = 98
$ws.Range("F49").Value = 425
$ws.Range("F50").Value = 47

# sheet2 (Worksheets item 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 83
$ws.Range("F10").Value = 276
$ws.Range("F22").Value = 3
$ws.Range("F24").Value = 6551

# sheet3 (Worksheets item 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F8").Value = 545
$ws.Range("F9").Value = 2126
$ws.Range("F11").Value = 979

# sheet4 (Worksheets item 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3302
$ws.Range("F5").Value = 7717
$ws.Range("F7").Value = 545
$ws.Range("F8").Value = 979
$ws.Range("F10").Value = 709
$ws.Range("F11").Value = 1128
$ws.Range("F12").Value = 1058
$ws.Range("F14").Value = 83
$ws.Range("F16").Value = 276
$ws.Range("F17").Value = 362
$ws.Range("F18").Value = 6179
$ws.Range("F19").Value = 76
$ws.Range("F23").Value = 5102
$ws.Range("F24").Value = 5925
$ws.Range("F25").Value = 360
$ws.Range("F26").Value = 160
$ws.Range("F30").Value = 8
$ws.Range("F31").Value = 1052
$ws.Range("F40").Value = 411
$ws.Range("F43").Value = 3
$ws.Range("F46").Value = 3202
$ws.Range("F47").Value = 98
$ws.Range("F48").Value = 47
